$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna5"
$ws.Range("C2").Value = "Epha5"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1261995
$ws.Range("H2").Value = 0.252399
$ws.Range("I2").Value = 0.07923641825693001
$ws.Range("J2").Value = 0.06507622895681928
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.006466
$ws.Range("N2").Value = 0.012932
$ws.Range("O2").Value = 0.008493427970384656
$ws.Range("P2").Value = 0.008493427970384656
$ws.Range("Q2").Value = 0.000816005967
$ws.Range("R2").Value = 0.003264023868
$ws.Range("S2").Value = 0.0006729888110965067
$ws.Range("T2").Value = 0.0005527202632290047

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efna5"
$ws.Range("C3").Value = "Epha5"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.5
$ws.Range("G3").Value = 0.1261995
$ws.Range("H3").Value = 0.252399
$ws.Range("I3").Value = 0.07923641825693001
$ws.Range("J3").Value = 0.06507622895681928
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.7548284999999999
$ws.Range("N3").Value = 1.509657
$ws.Range("O3").Value = 0.9915065720296153
$ws.Range("P3").Value = 0.9915065720296153
$ws.Range("Q3").Value = 0.09525897928574999
$ws.Range("R3").Value = 0.3810359171429999
$ws.Range("S3").Value = 0.0785634294458335
$ws.Range("T3").Value = 0.06452350869359028

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Efna5"
$ws.Range("C4").Value = "Epha5"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.6931216666666667
$ws.Range("H4").Value = 2.079365
$ws.Range("I4").Value = 0.4351877644756155
$ws.Range("J4").Value = 0.5361242826825642
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = 0.006466
$ws.Range("N4").Value = 0.012932
$ws.Range("O4").Value = 0.008493427970384656
$ws.Range("P4").Value = 0.008493427970384656
$ws.Range("Q4").Value = 0.004481724696666667
$ws.Range("R4").Value = 0.02689034818
$ws.Range("S4").Value = 0.003696235931166362
$ws.Range("T4").Value = 0.004553532978138501

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efna5"
$ws.Range("C5").Value = "Epha5"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.6931216666666667
$ws.Range("H5").Value = 2.079365
$ws.Range("I5").Value = 0.4351877644756155
$ws.Range("J5").Value = 0.5361242826825642
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.7548284999999999
$ws.Range("N5").Value = 1.509657
$ws.Range("O5").Value = 0.9915065720296153
$ws.Range("P5").Value = 0.9915065720296153
$ws.Range("Q5").Value = 0.5231879879675
$ws.Range("R5").Value = 3.139127927805
$ws.Range("S5").Value = 0.4314915285444491
$ws.Range("T5").Value = 0.5315707497044256

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Efna5"
$ws.Range("C6").Value = "Epha5"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.7733745
$ws.Range("H6").Value = 1.546749
$ws.Range("I6").Value = 0.4855758172674545
$ws.Range("J6").Value = 0.3987994883606166
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.5
$ws.Range("M6").Value = 0.006466
$ws.Range("N6").Value = 0.012932
$ws.Range("O6").Value = 0.008493427970384656
$ws.Range("P6").Value = 0.008493427970384656
$ws.Range("Q6").Value = 0.005000639517
$ws.Range("R6").Value = 0.020002558068
$ws.Range("S6").Value = 0.004124203228121786
$ws.Range("T6").Value = 0.003387174729017151

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Efna5"
$ws.Range("C7").Value = "Epha5"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.7733745
$ws.Range("H7").Value = 1.546749
$ws.Range("I7").Value = 0.4855758172674545
$ws.Range("J7").Value = 0.3987994883606166
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.7548284999999999
$ws.Range("N7").Value = 1.509657
$ws.Range("O7").Value = 0.9915065720296153
$ws.Range("P7").Value = 0.9915065720296153
$ws.Range("Q7").Value = 0.5837651137732499
$ws.Range("R7").Value = 2.335060455093
$ws.Range("S7").Value = 0.4814516140393327
$ws.Range("T7").Value = 0.3954123136315994
